$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" '64.004.96'
Set-TextValue "E2" '  +2.08%  '
Set-TextValue "D3" '3.417.32'
Set-TextValue "E3" '  +1.99%  '
Set-TextValue "E4" '  +0.00%  '
Set-TextValue "D5" '570.56'
Set-TextValue "D6" '157.04'
Set-TextValue "E6" '  +3.12%  '
Set-TextValue "D7" '1.00'
Set-TextValue "E7" '  +0.01%  '
Set-TextValue "D8" '3.415.14'
Set-TextValue "E8" '  +1.83%  '
Set-TextValue "D9" '0.551'
Set-TextValue "E9" '  +4.76%  '
Set-TextValue "E10" '  -1.26%  '
Set-TextValue "E11" '  +4.46%  '
Set-TextValue "E12" '  -1.30%  '
Set-TextValue "D13" '4.000.76'
Set-TextValue "E13" '  +2.02%  '
Set-TextValue "E14" '  -3.47%  '
Set-TextValue "E15" '  +6.73%  '
Set-TextValue "D16" '27.53'
Set-TextValue "E16" '  +2.59%  '
Set-TextValue "D17" '63.941.40'
Set-TextValue "E17" '  +1.97%  '
Set-TextValue "D18" '3.354.63'
Set-TextValue "E18" '  -0.46%  '
Set-TextValue "D19" '6.29'
Set-TextValue "E19" '  -0.69%  '
Set-TextValue "D20" '14.07'
Set-TextValue "D21" '382.73'
Set-TextValue "E21" '  -0.19%  '
Set-TextValue "D22" '8.08'
Set-TextValue "E22" '  -3.64%  '
Set-TextValue "D23" '0.999'
Set-TextValue "E23" '  -0.05%  '
Set-TextValue "D24" '72.15'
Set-TextValue "E24" '  +2.99%  '
Set-TextValue "E25" '  -0.21%  '
Set-TextValue "E26" '  +24.41%  '
Set-TextValue "D27" '9.43'
Set-TextValue "E27" '  +1.01%  '
Set-TextValue "D28" '0.179'
Set-TextValue "E28" '  -0.19%  '
Set-TextValue "E29" '  +0.15%  '
Set-TextValue "D30" '6.19'
Set-TextValue "E30" '  +10.80%  '
Set-TextValue "D31" '1.37'
Set-TextValue "E31" '  +5.23%  '
Set-TextValue "D32" '2.01'
Set-TextValue "E32" '  +0.20%  '
Set-TextValue "B33" 'RenderToken'
Set-TextValue "C33" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D33" '6.45'
Set-TextValue "E33" '  +1.84%  '
Set-TextValue "B34" 'EthereumClassic'
Set-TextValue "C34" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D34" '23.29'
Set-TextValue "E34" '  +1.28%  '
Set-TextValue "D35" '0.999'
Set-TextValue "E35" '  -0.04%  '
Set-TextValue "D36" '6.82'
Set-TextValue "E36" '  +1.59%  '
Set-TextValue "D37" '159.45'
Set-TextValue "E37" '  +1.05%  '
Set-TextValue "E38" '  +0.14%  '
Set-TextValue "D39" '2.970.15'
Set-TextValue "E39" '  +7.00%  '
Set-TextValue "E40" '  +2.92%  '
Set-TextValue "D41" '1.84'
Set-TextValue "E41" '  -1.95%  '
Set-TextValue "D42" '26.86'
Set-TextValue "E42" '  -0.57%  '
Set-TextValue "D43" '42.19'
Set-TextValue "E43" '  +3.14%  '
Set-TextValue "E44" '  -4.47%  '
Set-TextValue "E45" '  +1.88%  '
Set-TextValue "E46" '  +1.05%  '
Set-TextValue "D47" '23.31'
Set-TextValue "E47" '  +6.44%  '
Set-TextValue "D48" '1.07'
Set-TextValue "E48" '  +2.61%  '
Set-TextValue "D49" '2.20'
Set-TextValue "E49" '  +21.44%  '
Set-TextValue "B50" 'Cosmos'
Set-TextValue "C50" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D50" '6.38'
Set-TextValue "E50" '  +0.68%  '
Set-TextValue "B51" 'SuiNetwork'
Set-TextValue "C51" 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue "D51" '0.835'
Set-TextValue "E51" '  +3.50%  '
